$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Date: 2025-05-21T14:22:51+00:00 -> 2025-06-13T15:45:04+00:00
$wsMeta.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# FHIR Version: 4.3.0 -> 4.0.1
$wsMeta.Range("B15").Value = "4.0.1"

# --- Elements sheet updates ---
$wsElem = $wb.Worksheets.Item("Elements")

# Extension ele-1 constraint (row 2, Constraint(s) column AJ):
# drop the "unless an empty Parameters resource ... or $this is Parameters" clause
$wsElem.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Extension.id Type(s) column (row 3, column K): id -> string
$wsElem.Range("K3").Value = "string`n"

# Extension.value[x] Definition column (row 6, column M): R4B -> R4 in URL
$wsElem.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
